$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Rushing" (sheet1) - Week 15 stat corrections / logging
# ---------------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# J.Fields (row 3)
$rushing.Range("C3").Value = 5
$rushing.Range("D3").Value = 7
$rushing.Range("E3").Value = 4
$rushing.Range("F3").Value = 4

# D.Montgomery (row 4)
$rushing.Range("C4").Value = 128
$rushing.Range("D4").Value = 64
$rushing.Range("E4").Value = 24
$rushing.Range("F4").Value = 22

# D.Mooney (row 9)
$rushing.Range("C9").Value = 3
$rushing.Range("F9").Value = 2

# J.Grant (row 11)
$rushing.Range("C11").Value = 3
$rushing.Range("E11").Value = 1

# ---------------------------------------------------------------------------
# Sheet "Receiving" (sheet2) - Week 15 stat corrections + simulated Week 16
# (new player row for J.Horsted)
# ---------------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# D.Montgomery (row 2)
$receiving.Range("C2").Value = 46
$receiving.Range("D2").Value = 40
$receiving.Range("G2").Value = 6
$receiving.Range("H2").Value = 6

# K.Herbert (row 4)
$receiving.Range("C4").Value = 8
$receiving.Range("D4").Value = 7

# D.Mooney (row 6)
$receiving.Range("C6").Value = 76
$receiving.Range("D6").Value = 44
$receiving.Range("E6").Value = 26
$receiving.Range("F6").Value = 13
$receiving.Range("G6").Value = 8
$receiving.Range("H6").Value = 4

# D.Byrd (row 8)
$receiving.Range("C8").Value = 22
$receiving.Range("D8").Value = 17
$receiving.Range("G8").Value = 3

# C.Kmet (row 10)
$receiving.Range("C10").Value = 68
$receiving.Range("D10").Value = 43
$receiving.Range("E10").Value = 10
$receiving.Range("F10").Value = 6
$receiving.Range("G10").Value = 10

# J.Graham (row 11)
$receiving.Range("C11").Value = 16
$receiving.Range("D11").Value = 8
$receiving.Range("G11").Value = 7

# New player J.Horsted (row 13) - copy row 12's look (border/bold on col A)
# then fill in the new values.
$receiving.Range("A12").Copy($receiving.Range("A13"))
$receiving.Range("A13").Value = 11
$receiving.Range("B13").Value = "J.Horsted"
$receiving.Range("C13").Value = 1
$receiving.Range("D13").Value = 1
$receiving.Range("E13").Value = 0
$receiving.Range("F13").Value = 0
$receiving.Range("G13").Value = 1
$receiving.Range("H13").Value = 1
